$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh price (D) and 1h volume change (E) figures for each coin row
$ws.Range("D2").Value = "26.992.87"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "1.557.28"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").Value = "'207.08"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("D9").Value = "'21.56"
$ws.Range("E9").Value = "  -0.24%  "

$ws.Range("D10").Value = "'0.0584"
$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "1.779.03"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "1.557.55"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("D14").Value = "'3.71"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D16").Value = "26.978.46"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").Value = "'61.73"
$ws.Range("E17").Value = "  +0.78%  "

$ws.Range("D18").Value = "'214.87"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").Value = "'7.26"
$ws.Range("E20").Value = "  -1.12%  "

$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("D22").Value = "'4.05"
$ws.Range("E22").Value = "  -1.79%  "

$ws.Range("D23").Value = "'9.21"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").Value = "'1.98"
$ws.Range("E24").Value = "  -2.07%  "

$ws.Range("D25").Value = "'153.29"
$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").Value = "'14.91"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").Value = "'1.01"
$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("E30").Value = "  -1.00%  "

$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").Value = "1.371.66"
$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("D34").Value = "'2.96"
$ws.Range("E34").Value = "  +1.45%  "

$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +2.66%  "

$ws.Range("D36").Value = "'0.967"
$ws.Range("E36").Value = "  +5.17%  "

$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("E38").Value = "  +0.76%  "

$ws.Range("D41").Value = "'1.01"
$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D44").Value = "'2.23"
$ws.Range("E44").Value = "  +2.68%  "

$ws.Range("D45").Value = "'63.84"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("E46").Value = "  -1.49%  "

$ws.Range("D47").Value = "1.692.15"
$ws.Range("E47").Value = "  -0.27%  "

$ws.Range("E48").Value = "  -3.56%  "

$ws.Range("D49").Value = "'86.34"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  +0.62%  "

# Coin ranking reordering: rows 39/40 and 42/43 swap coins, with updated data
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.520"
$ws.Range("E39").Value = "  -1.34%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.810"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.51"
$ws.Range("E42").Value = "  -0.34%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'0.980"
$ws.Range("E43").Value = "  -1.48%  "
